# Updated symbol list on Mon Dec 26 16:53:09 UTC 2022 with GitHub Actions
#
# Refreshes the coinranking snapshot: new "Price" (column D) readings and a
# couple of corrected "Volume(1h)" (column E) labels.
#
# Column D is stored as TEXT (coinranking prices are captured verbatim, so
# trailing zeros / significant digits must survive exactly, e.g. "0.8100"
# must stay "0.8100" and not collapse to the number 0.81). Plainly assigning
# a numeric-looking string to `.Value` lets Excel auto-coerce it to a
# Number, which is wrong here, so each cell is explicitly forced to Text
# (NumberFormat "@") before the new value is written, then its style is put
# back to "Normal" so we don't leave a stray number-format behind on the
# cell.
#
# Column E holds plain (non-numeric) text, so it can be assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextPrice($a1, $value) {
    $cell = $ws.Range($a1)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextPrice "D2"  "242.61"
Set-TextPrice "D4"  "5.426"
Set-TextPrice "D5"  "0.05898"
Set-TextPrice "D6"  "3.436"
Set-TextPrice "D7"  "6.531"
Set-TextPrice "D8"  "0.8100"
Set-TextPrice "D9"  "0.9730"
Set-TextPrice "D11" "0.07402"
Set-TextPrice "D12" "0.03277"
Set-TextPrice "D13" "0.03053"
Set-TextPrice "D14" "0.09340"
Set-TextPrice "D15" "3.859"
Set-TextPrice "D16" "0.001578"
Set-TextPrice "D17" "0.04669"
Set-TextPrice "D18" "0.0005930"
Set-TextPrice "D19" "0.005847"
Set-TextPrice "D20" "0.001267"
Set-TextPrice "D21" "0.004899"
Set-TextPrice "D22" "0.00006802"
Set-TextPrice "D24" "2.170"

Set-TextPrice "D27" "0.0002284"
$ws.Range("E27").Value = "26UpBotsUBXTWorstin24h"

Set-TextPrice "D40" "0.03932"
Set-TextPrice "D41" "0.006189"
Set-TextPrice "D42" "0.1072"
Set-TextPrice "D43" "0.002561"
Set-TextPrice "D45" "0.00005175"

Set-TextPrice "D47" "0.6660"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

Set-TextPrice "D48" "0.002382"
Set-TextPrice "D49" "0.00002100"
Set-TextPrice "D50" "0.0002000"
